$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.351.42'
$ws.Range('E2').Value = '  +0.95%  '

$ws.Range('D3').Value = '3.922.46'
$ws.Range('E3').Value = '  -0.74%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '485.36'
$ws.Range('E5').Value = '  -0.03%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.34'
$ws.Range('E6').Value = '  -1.55%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.625'
$ws.Range('E7').Value = '  -0.38%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.997'
$ws.Range('E8').Value = '  -0.12%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.735'
$ws.Range('E9').Value = '  +0.49%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.168'
$ws.Range('E10').Value = '  -2.49%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0000349'
$ws.Range('E11').Value = '  -6.12%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '43.19'
$ws.Range('E12').Value = '  -0.66%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '10.80'
$ws.Range('E13').Value = '  +3.52%  '

$ws.Range('D14').Value = '4.543.39'
$ws.Range('E14').Value = '  -1.17%  '

$ws.Range('D15').Value = '3.932.39'
$ws.Range('E15').Value = '  +0.39%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.37'
$ws.Range('E16').Value = '  -3.63%  '

$ws.Range('E17').Value = '  -0.78%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '20.13'
$ws.Range('E18').Value = '  +1.07%  '

$ws.Range('E19').Value = '  -0.26%  '

$ws.Range('D20').Value = '68.407.85'
$ws.Range('E20').Value = '  +0.81%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '432.76'
$ws.Range('E21').Value = '  -0.12%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.52'
$ws.Range('E22').Value = '  +3.22%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '15.13'
$ws.Range('E23').Value = '  +5.01%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '88.46'
$ws.Range('E24').Value = '  +0.53%  '

$ws.Range('E25').Value = '  +22.73%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.24'
$ws.Range('E26').Value = '  +12.08%  '

$ws.Range('E27').Value = '  +1.45%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '37.97'
$ws.Range('E28').Value = '  -1.76%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.70'
$ws.Range('E29').Value = '  -0.90%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '720.48'
$ws.Range('E30').Value = '  -0.13%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '13.85'
$ws.Range('E31').Value = '  +3.98%  '

$ws.Range('E32').Value = '  +0.62%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.93'
$ws.Range('E33').Value = '  +3.78%  '

$ws.Range('B34').Value = 'PEPE'
$ws.Range('C34').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D34').Value = '0.0₃0924'
$ws.Range('E34').Value = '  +2.08%  '

$ws.Range('B35').Value = 'NEARProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.15'
$ws.Range('E35').Value = '  +15.03%  '

$ws.Range('B36').Value = 'InjectiveProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '41.77'
$ws.Range('E36').Value = '  -0.36%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '61.13'
$ws.Range('E37').Value = '  +1.22%  '

$ws.Range('E38').Value = '  +21.84%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.398'
$ws.Range('E39').Value = '  +18.94%  '

$ws.Range('B40').Value = 'Kaspa'
$ws.Range('C40').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.146'
$ws.Range('E40').Value = '  -3.94%  '

$ws.Range('B41').Value = 'Dai'
$ws.Range('C41').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.999'
$ws.Range('E41').Value = '  +0.08%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0493'
$ws.Range('E42').Value = '  +3.87%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.14'
$ws.Range('E43').Value = '  +3.20%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.01'
$ws.Range('E44').Value = '  +4.68%  '

$ws.Range('E45').Value = '  +0.32%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.36'
$ws.Range('E46').Value = '  +4.23%  '

$ws.Range('E47').Value = '  -0.09%  '

$ws.Range('E48').Value = '  -0.17%  '

$ws.Range('E49').Value = '  -4.00%  '

$ws.Range('B50').Value = 'Monero'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '145.12'
$ws.Range('E50').Value = '  -2.49%  '

$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').Value = '0.0₆0341'
$ws.Range('E51').Value = '  +25.23%  '
